# Refresh cryptos list snapshot (price + 1h volume change) on Sheet1.
# Values that could otherwise be auto-parsed as numbers by Excel are written
# with a leading apostrophe to force text, then ClearFormats() removes the
# resulting quote-prefix style so the cell keeps its original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.633.96"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.596.36"
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'211.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").Value = "'0.514"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.02%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").Value = "'19.53"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.24%  "
$ws.Range("D11").Value = "'0.0837"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").Value = "1.820.90"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.605.84"
$ws.Range("E13").Value = "  +1.09%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'4.03"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "'64.40"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "26.614.98"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "'208.78"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").Value = "'6.94"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.56%  "
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("E23").Value = "  -2.49%  "
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").Value = "'145.18"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").Value = "'7.15"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").Value = "'15.26"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'0.0507"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "'0.661"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.98%  "
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").Value = "1.277.98"
$ws.Range("D36").Value = "'2.44"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("E37").Value = "  +0.64%  "
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").Value = "'0.844"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("D41").Value = "'5.48"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("E42").Value = "  +1.76%  "
$ws.Range("D43").Value = "'64.50"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("D44").Value = "'0.785"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("D45").Value = "'0.919"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +9.36%  "
$ws.Range("D46").Value = "1.733.49"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").Value = "'89.96"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  -2.33%  "
$ws.Range("E50").Value = "  +4.55%  "
$ws.Range("E51").Value = "  +0.39%  "

Write-Host "Applied cryptos update."